# [Feat] NPC 대화 조정
#
# BlackSmith.xlsx - Sheet1 dialogue table.
# Replace the throwaway test lines in column A (rows 1, 2 and 4) with the
# real NPC greeting lines; row 3 keeps its existing text. This mirrors the
# shared-string content/order changes in the target OOXML diff:
#   - drop the two "test" strings ("실험용 1회용 대사" / "오늘 날씨는 어때")
#   - rename the old "반복용 대사 없는거 빼곤 다 있다고" line to "오늘도 왔는가"
#   - append two new repeat-dialogue lines
#
# Also move the active selection to B8, matching the saved sheet view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value2 = "오늘도 왔는가"
$ws.Cells.Item(2, 1).Value2 = "요즘 농사는 괜찮게 되고 있지?"
$ws.Cells.Item(4, 1).Value2 = "우리 마을의 유일 마트라고"

$ws.Range("B8").Select()
